$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "O-O": only the "N. atoms" text (column D) formatting changes -
# each line now ends with "; " before the newline / end of string.
# ---------------------------------------------------------------------------
$wsOO = $wb.Worksheets.Item("O-O")

$wsOO.Range("D2").Value = "12 atoms have CN of 0; `n32 atoms have CN of 1; "
$wsOO.Range("D3").Value = "4 atoms have CN of 2; `n18 atoms have CN of 4; `n16 atoms have CN of 6; `n6 atoms have CN of 8; "
$wsOO.Range("D4").Value = "32 atoms have CN of 1; `n12 atoms have CN of 2; "
$wsOO.Range("D5").Value = "14 atoms have CN of 0; `n20 atoms have CN of 1; `n10 atoms have CN of 2; "
$wsOO.Range("D6").Value = "28 atoms have CN of 1; `n10 atoms have CN of 2; `n6 atoms have CN of 3; "
$wsOO.Range("D7").Value = "14 atoms have CN of 2; `n8 atoms have CN of 3; `n20 atoms have CN of 4; `n2 atoms have CN of 6; "

# ---------------------------------------------------------------------------
# Sheet "Pb-Pb": updated average CN values (column C) and shell-fitting
# breakdown text (column D).
# ---------------------------------------------------------------------------
$wsPbPb = $wb.Worksheets.Item("Pb-Pb")

$wsPbPb.Range("C2").Value = 0.8
$wsPbPb.Range("D2").Value = "4 atoms have CN of 0; `n10 atoms have CN of 1; `n1 atoms have CN of 2; "

$wsPbPb.Range("C3").Value = 4.266666666666667
$wsPbPb.Range("D3").Value = "14 atoms have CN of 4; `n1 atoms have CN of 8; "

$wsPbPb.Range("C4").Value = 1.6
$wsPbPb.Range("D4").Value = "2 atoms have CN of 0; `n4 atoms have CN of 1; `n8 atoms have CN of 2; `n1 atoms have CN of 4; "

# ---------------------------------------------------------------------------
# Sheet "Pb-O": updated average CN values (column C) and shell-fitting
# breakdown text (column D).
# ---------------------------------------------------------------------------
$wsPbO = $wb.Worksheets.Item("Pb-O")

$wsPbO.Range("C2").Value = 5.2
$wsPbO.Range("D2").Value = "4 atoms have CN of 3; `n11 atoms have CN of 6; "

$wsPbO.Range("C3").Value = 2.4
$wsPbO.Range("D3").Value = "12 atoms have CN of 2; `n3 atoms have CN of 4; "

$wsPbO.Range("C4").Value = 2.4
$wsPbO.Range("D4").Value = "10 atoms have CN of 2; `n4 atoms have CN of 3; `n1 atoms have CN of 4; "

$wsPbO.Range("C5").Value = 4.266666666666667
$wsPbO.Range("D5").Value = "14 atoms have CN of 4; `n1 atoms have CN of 8; "

$wsPbO.Range("C6").Value = 0.9333333333333333
$wsPbO.Range("D6").Value = "4 atoms have CN of 0; `n8 atoms have CN of 1; `n3 atoms have CN of 2; "

# ---------------------------------------------------------------------------
# Sheet "O-Pb": updated average CN values (column C) and shell-fitting
# breakdown text (column D).
# ---------------------------------------------------------------------------
$wsOPb = $wb.Worksheets.Item("O-Pb")

$wsOPb.Range("C2").Value = 1.772727272727273
$wsOPb.Range("D2").Value = "20 atoms have CN of 1; `n14 atoms have CN of 2; `n10 atoms have CN of 3; "

$wsOPb.Range("C3").Value = 0.8181818181818182
$wsOPb.Range("D3").Value = "16 atoms have CN of 0; `n20 atoms have CN of 1; `n8 atoms have CN of 2; "

$wsOPb.Range("C4").Value = 0.8181818181818182
$wsOPb.Range("D4").Value = "10 atoms have CN of 0; `n32 atoms have CN of 1; `n2 atoms have CN of 2; "

$wsOPb.Range("C5").Value = 1.454545454545455
$wsOPb.Range("D5").Value = "10 atoms have CN of 0; `n8 atoms have CN of 1; `n24 atoms have CN of 2; `n2 atoms have CN of 4; "

$wsOPb.Range("C6").Value = 0.3181818181818182
$wsOPb.Range("D6").Value = "30 atoms have CN of 0; `n14 atoms have CN of 1; "
